# Build v2.1.2: Fix SearchCriteria variants and Schemas sheet grouping/sorting
#
# This script rewrites the "schema summary" row (row 3) on the Body sheet and
# on each HTTP response-code sheet so that it references the correct request /
# response / error schema name, and removes the now-redundant detail rows
# that used to spell out the schema's individual fields (those now live only
# inside the referenced schema definition).

$wb = $excel.ActiveWorkbook

function Set-SchemaRow {
    param(
        [string]$SheetName,
        [string]$SchemaName,
        [string]$SectionLabel
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Drop every data row below the schema-summary row (row 3) - the detailed
    # field-by-field rows are superseded by the schema reference itself.
    $ws.Range("A4:O1048576").Clear()

    # Column A is the Section bucket label ("body" on the request sheet,
    # "content" on every response sheet).
    $ws.Cells.Item(3, 1).Value = $SectionLabel                # A: Section
    $ws.Cells.Item(3, 2).Value = $SchemaName                # B: Name
    $ws.Cells.Item(3, 3).ClearContents()                     # C: Parent
    $ws.Cells.Item(3, 4).ClearContents()                     # D: Description
    $ws.Cells.Item(3, 5).Value = "schema"                    # E: Type
    $ws.Cells.Item(3, 6).ClearContents()                     # F: Items Data Type
    $ws.Cells.Item(3, 7).Value = $SchemaName                 # G: Schema Name
    $ws.Cells.Item(3, 8).ClearContents()                     # H: Format
    $ws.Cells.Item(3, 9).Value = "Yes"                       # I: Mandatory
    $ws.Cells.Item(3, 10).ClearContents()                    # J: Min
    $ws.Cells.Item(3, 11).ClearContents()                    # K: Max
    $ws.Cells.Item(3, 12).ClearContents()                    # L: PatternEba
    $ws.Cells.Item(3, 13).ClearContents()                    # M: Regex
    $ws.Cells.Item(3, 14).ClearContents()                    # N: Allowed value
    $ws.Cells.Item(3, 15).ClearContents()                    # O: Example
}

# Request body -> references the request schema, drops the old settlementBIC row.
Set-SchemaRow "Body" "getDefaultAgenda.211207Request" "body"

# Successful responses -> reference the response schema.
Set-SchemaRow "200" "getDefaultAgenda.211207Response" "content"
Set-SchemaRow "204" "getDefaultAgenda.211207Response" "content"

# Canonical error payload definition.
Set-SchemaRow "400" "errorResponse" "content"

# Remaining error responses reference the shared error-response variant schema.
Set-SchemaRow "401" "errorResponse1" "content"
Set-SchemaRow "403" "errorResponse1" "content"
Set-SchemaRow "404" "errorResponse1" "content"
Set-SchemaRow "429" "errorResponse1" "content"
Set-SchemaRow "500" "errorResponse1" "content"
